$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 107 (pushes old rows 107-200 down to 108-201)
$ws.Rows("107:107").Insert()

# Populate the newly inserted row 107 with the new data entry
$ws.Range("A107").Value = 10
$ws.Range("B107").Value = "Vega Modelo de Temuco"
$ws.Range("C107").Value = "La Araucanía"
$ws.Range("D107").Value = 45062
$ws.Range("E107").Value = 9
$ws.Range("F107").Value = 100114002
$ws.Range("G107").Value = "Camote"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 5
$ws.Range("K107").Value = 26000
$ws.Range("L107").Value = 26000
$ws.Range("M107").Value = 26000
$ws.Range("N107").Value = "$/caja 18 kilos"
$ws.Range("O107").Value = "Perú"
$ws.Range("P107").Value = 1444
$ws.Range("Q107").Value = 18
$ws.Range("R107").Value = "Hortaliza"
